# Lotto 964~975 tracking sheet update ("Add files via upload")
# - Fill in the previously-blank "결과" (result) row for 964회차 (row 14, T:Z)
#   and the previously-blank 971회차 block (rows 16-20, columns B:G).
# - Highlight the matched numbers in the 971회차 "올" row (row 12: U12, V12, X12)
#   with the same bold + yellow-fill look used elsewhere on the sheet for hits.
# - Move the active selection to E21, matching where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 964회차 block (row 14): fill in the "결과"(result)+"보너스"(bonus) row, S:Z
# ---------------------------------------------------------------------------
$ws.Range("T14").Value = 3
$ws.Range("U14").Value = 9
$ws.Range("V14").Value = 10
$ws.Range("W14").Value = 29
$ws.Range("X14").Value = 40
$ws.Range("Y14").Value = 45
$ws.Range("Z14").Value = 7

# ---------------------------------------------------------------------------
# 971회차 block: fill in the previously-empty "올"(row16)/A(row17)/
# "랜덤"(row18)/A(row19)/"결과"(row20) number grid, columns B:G
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 26
$ws.Range("E16").Value = 31
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 44

$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 29
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 44

$ws.Range("B18").Value = 8
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = 43
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 39

$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 43
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 45

$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = 27
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 42

# ---------------------------------------------------------------------------
# Highlight the numbers in row 12 ("올") that hit the 964회차 result (row 14):
# U12=9, V12=10 and X12=45 all now match, so bold them with the sheet's
# yellow "hit" fill (same treatment as the other highlighted hit-cells).
# ---------------------------------------------------------------------------
foreach ($addr in @("U12", "V12", "X12")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Interior.Color = 65535
}

# ---------------------------------------------------------------------------
# Leave the selection where the editor last left it.
# ---------------------------------------------------------------------------
[void]$ws.Range("E21").Select()
